# Auxiliary output from ELC_DMD process added with corresponding changes
# across documents (Scenario EXP file and SubRes DK price).
#
# Logical changes in this workbook (SubRes_ELC_DKprices.xlsx):
#   - ELC_IMP_EXP_ProcsR!D3: "Export of electricity to Denmark West"
#         -> "Sale of electricity to Denmark West"
#   - ELC_IMP_EXP_TechsR!C3: "Export of electricity to Denmark West"
#         -> "Sale of electricity to Denmark West"
#   - ELC_IMP_EXP_TechsR!D3: "ELC_DEM" -> "ELC_PRIS"
#   - Selection/active-cell bookkeeping updated on both sheets.

$wb = $excel.ActiveWorkbook

$wsProcs = $wb.Worksheets.Item("ELC_IMP_EXP_ProcsR")
$wsTechs = $wb.Worksheets.Item("ELC_IMP_EXP_TechsR")

# --- ELC_IMP_EXP_ProcsR: rename the EXP process description ---
$wsProcs.Range("D3").Value = "Sale of electricity to Denmark West"

# --- ELC_IMP_EXP_TechsR: rename the EXPELC-DKW tech description and its
#     output commodity (from electricity demand to electricity price) ---
$wsTechs.Range("C3").Value = "Sale of electricity to Denmark West"
$wsTechs.Range("D3").Value = "ELC_PRIS"

# --- Update the last active selection on each sheet to match the edit ---
$wsTechs.Activate()
$wsTechs.Range("D3").Select()

$wsProcs.Activate()
$wsProcs.Range("G13").Select()
